$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'317.82"
$ws.Range("E2").Value = "'3.77%"
$ws.Range("D3").Value = "'39.69"
$ws.Range("E3").Value = "'1.85%"
$ws.Range("D4").Value = "'5.135"
$ws.Range("E4").Value = "'0.86%"
$ws.Range("D5").Value = "'0.08209"
$ws.Range("E5").Value = "'1.88%"
$ws.Range("D6").Value = "'2.018"
$ws.Range("E6").Value = "'4.98%"
$ws.Range("D7").Value = "'8.305"
$ws.Range("E7").Value = "'4.27%"
$ws.Range("D8").Value = "'0.9316"
$ws.Range("E8").Value = "'-0.24%"
$ws.Range("E9").Value = "'-2.36%"
$ws.Range("D10").Value = "'0.1994"
$ws.Range("E10").Value = "'3.57%"
$ws.Range("D11").Value = "'0.09079"
$ws.Range("E11").Value = "'-0.17%"
$ws.Range("D12").Value = "'0.03476"
$ws.Range("E12").Value = "'-0.95%"
$ws.Range("D13").Value = "'0.09805"
$ws.Range("E13").Value = "'0.13%"
$ws.Range("D14").Value = "'0.001393"
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("D15").Value = "'0.006203"
$ws.Range("E15").Value = "'4.72%"
$ws.Range("D16").Value = "'3.686"
$ws.Range("E16").Value = "'-2.88%"
$ws.Range("D17").Value = "'4.285"
$ws.Range("E17").Value = "'2.10%"
$ws.Range("D18").Value = "'3.312"
$ws.Range("E18").Value = "'-3.01%"
$ws.Range("D19").Value = "'0.3474"
$ws.Range("E19").Value = "'0.89%"
$ws.Range("E20").Value = "'-2.68%"
$ws.Range("D21").Value = "'4.897"
$ws.Range("E21").Value = "'2.24%"
$ws.Range("D22").Value = "'0.2449"
$ws.Range("E22").Value = "'-2.45%"
$ws.Range("D23").Value = "'0.04328"
$ws.Range("E23").Value = "'-1.16%"
$ws.Range("D24").Value = "'0.001225"
$ws.Range("E24").Value = "'-1.15%"
$ws.Range("D25").Value = "'0.004770"
$ws.Range("E25").Value = "'11.71%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'-0.34%"
$ws.Range("D27").Value = "'0.0003997"
$ws.Range("E27").Value = "'-10.14%"
$ws.Range("D39").Value = "'0.02213"
$ws.Range("E39").Value = "'8.37%"
$ws.Range("D40").Value = "'0.05234"
$ws.Range("E40").Value = "'3.92%"
$ws.Range("D41").Value = "'0.007508"
$ws.Range("E41").Value = "'0.98%"
$ws.Range("D42").Value = "'0.009664"
$ws.Range("E42").Value = "'-4.53%"
$ws.Range("E43").Value = "'2.28%"
$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'-1.50%"
$ws.Range("D45").Value = "'0.009851"
$ws.Range("E45").Value = "'8.55%"
$ws.Range("D46").Value = "'0.00006587"
$ws.Range("E46").Value = "'6.22%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.33%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.001200"
$ws.Range("E48").Value = "'-25.07%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.002763"
$ws.Range("E49").Value = "'-1.53%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.33%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.33%"
